# Re-orders the data rows (2-45) of the "Artfynd" sheet: the content that
# used to live in one row is moved to another row, while the row index
# (and therefore any row-level formatting) stays put. Columns Y and AA
# hold textual dates ("2021-10-11"); Excel auto-coerces such strings to
# real dates when they're pushed back through Value2, so those two
# columns are pre-formatted as Text to keep them as plain strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 45
$dataRange = $ws.Range("A$firstRow`:AY$lastRow")

# sourceOffsets[i] (0-based) = 0-based row offset (relative to $firstRow)
# whose data should end up at target offset i (relative to $firstRow).
$sourceOffsets = @(6,31,32,33,34,7,35,36,8,0,9,10,1,2,3,37,38,4,39,40,41,42,11,5,43,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30)

# Keep the Startdatum/Slutdatum columns as Text so Excel doesn't turn the
# "yyyy-mm-dd" strings into real date serials when we write them back.
$ws.Range("Y$firstRow`:Y$lastRow").NumberFormat = "@"
$ws.Range("AA$firstRow`:AA$lastRow").NumberFormat = "@"

$values = $dataRange.Value2
$rowCount = $values.GetLength(0)
$colCount = $values.GetLength(1)

# $values (Value2's result) uses 1-based COM-style indexing; a freshly
# New-Object'd .NET array uses 0-based indexing, and Excel's interop maps
# that straight onto the target range (row 0 -> first row of the range).
$newValues = New-Object 'object[,]' $rowCount, $colCount
for ($i = 1; $i -le $rowCount; $i++) {
    $srcRow = $sourceOffsets[$i - 1] + 1
    for ($j = 1; $j -le $colCount; $j++) {
        $newValues[$i - 1, $j - 1] = $values[$srcRow, $j]
    }
}

$dataRange.Value2 = $newValues
